# Scheduled-runner data refresh: re-pulls Universalis market-board averages
# for the Leve profit tables and rewrites the dependent price/profit columns
# (H:currentAveragePrice, I/J: NQ/HQ average, K/L: Leve price NQ/HQ,
# M/N: Leve profit NQ/HQ) on every crafting-job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 10010.5
$ws.Range("I20").Value = 10010.5
$ws.Range("K20").Value = 10010.5
$ws.Range("M20").Value = -9780.5
# Row 35
$ws.Range("H35").Value = 10010.5
$ws.Range("I35").Value = 10010.5
$ws.Range("K35").Value = 10010.5
$ws.Range("M35").Value = -9631.5
# Row 135
$ws.Range("H135").Value = 2790.106
$ws.Range("I135").Value = 2717.2407
$ws.Range("J135").Value = 3118
$ws.Range("K135").Value = 24455.1663
$ws.Range("L135").Value = 28062
$ws.Range("M135").Value = -21920.1663
$ws.Range("N135").Value = -33132
# Row 137
$ws.Range("H137").Value = 213794.05
$ws.Range("I137").Value = 291565.6
$ws.Range("J137").Value = 1218.5333
$ws.Range("K137").Value = 874696.7999999999
$ws.Range("L137").Value = 3655.5999
$ws.Range("M137").Value = -872146.7999999999
$ws.Range("N137").Value = -8755.599900000001
# Row 138
$ws.Range("H138").Value = 1330.45
$ws.Range("I138").Value = 725.4211
$ws.Range("J138").Value = 2132.465
$ws.Range("K138").Value = 2176.2633
$ws.Range("L138").Value = 6397.395
$ws.Range("M138").Value = 2963.7367
$ws.Range("N138").Value = -16677.395

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3014.4146
$ws.Range("I61").Value = 3154.3513
$ws.Range("J61").Value = 1720
$ws.Range("K61").Value = 3154.3513
$ws.Range("L61").Value = 1720
$ws.Range("M61").Value = -2942.3513
$ws.Range("N61").Value = -2144
# Row 102
$ws.Range("H102").Value = 117648280
$ws.Range("I102").Value = 181819330
$ws.Range("K102").Value = 181819330
$ws.Range("M102").Value = -181817708
# Row 110
$ws.Range("H110").Value = 719.61536
$ws.Range("I110").Value = 671.2778
$ws.Range("J110").Value = 828.375
$ws.Range("K110").Value = 671.2778
$ws.Range("L110").Value = 828.375
$ws.Range("M110").Value = 1373.7222
$ws.Range("N110").Value = -4918.375
# Row 132
$ws.Range("H132").Value = 5954399.5
$ws.Range("I132").Value = 6946132
$ws.Range("K132").Value = 20838396
$ws.Range("M132").Value = -20835866
# Row 136
$ws.Range("H136").Value = 3014.4146
$ws.Range("I136").Value = 3154.3513
$ws.Range("J136").Value = 1720
$ws.Range("K136").Value = 9463.053899999999
$ws.Range("L136").Value = 5160
$ws.Range("M136").Value = -6913.053899999999
$ws.Range("N136").Value = -10260

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1119.7188
$ws.Range("I105").Value = 1013.34784
$ws.Range("J105").Value = 1391.5555
$ws.Range("K105").Value = 1013.34784
$ws.Range("L105").Value = 1391.5555
$ws.Range("M105").Value = 733.65216
$ws.Range("N105").Value = -4885.5555

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 19237976
$ws.Range("I31").Value = 26316930
$ws.Range("J31").Value = 23671.143
$ws.Range("K31").Value = 26316930
$ws.Range("L31").Value = 23671.143
$ws.Range("M31").Value = -26316635
$ws.Range("N31").Value = -24261.143
# Row 34
$ws.Range("H34").Value = 19237976
$ws.Range("I34").Value = 26316930
$ws.Range("J34").Value = 23671.143
$ws.Range("K34").Value = 26316930
$ws.Range("L34").Value = 23671.143
$ws.Range("M34").Value = -26316728
$ws.Range("N34").Value = -24075.143
# Row 52
$ws.Range("H52").Value = 38880
$ws.Range("J52").Value = 38880
$ws.Range("L52").Value = 38880
$ws.Range("N52").Value = -39468
# Row 58
$ws.Range("H58").Value = 3692091.2
$ws.Range("I58").Value = 4231963.5
$ws.Range("K58").Value = 4231963.5
$ws.Range("M58").Value = -4231760.5
# Row 132
$ws.Range("H132").Value = 7096136
$ws.Range("I132").Value = 9010688
$ws.Range("K132").Value = 27032064
$ws.Range("M132").Value = -27029534
# Row 134
$ws.Range("H134").Value = 20193248
$ws.Range("I134").Value = 23585788
$ws.Range("J134").Value = 5209534.5
$ws.Range("K134").Value = 70757364
$ws.Range("L134").Value = 15628603.5
$ws.Range("M134").Value = -70754829
$ws.Range("N134").Value = -15633673.5
# Row 136
$ws.Range("H136").Value = 3692091.2
$ws.Range("I136").Value = 4231963.5
$ws.Range("K136").Value = 12695890.5
$ws.Range("M136").Value = -12693340.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 615.86664
$ws.Range("I5").Value = 432
$ws.Range("K5").Value = 1296
$ws.Range("M5").Value = -1184
# Row 122
$ws.Range("H122").Value = 727
$ws.Range("J122").Value = 733.1111
$ws.Range("L122").Value = 6597.9999
$ws.Range("N122").Value = -11497.9999
# Row 131
$ws.Range("H131").Value = 13213348
$ws.Range("I131").Value = 45455040
$ws.Range("J131").Value = 1391393.9
$ws.Range("K131").Value = 136365120
$ws.Range("L131").Value = 4174181.7
$ws.Range("M131").Value = -136360080
$ws.Range("N131").Value = -4184261.7
# Row 135
$ws.Range("H135").Value = 615.86664
$ws.Range("I135").Value = 432
$ws.Range("K135").Value = 3888
$ws.Range("M135").Value = -1353

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1219.3
$ws.Range("I113").Value = 1176.2142
$ws.Range("K113").Value = 1176.2142
$ws.Range("M113").Value = 993.7858000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1469.9286
$ws.Range("I61").Value = 1051.7273
$ws.Range("K61").Value = 1051.7273
$ws.Range("M61").Value = -849.7273
# Row 113
$ws.Range("H113").Value = 1469.9286
$ws.Range("I113").Value = 1051.7273
$ws.Range("K113").Value = 1051.7273
$ws.Range("M113").Value = 1118.2727
# Row 132
$ws.Range("H132").Value = 2818008.8
$ws.Range("I132").Value = 3449039.5
$ws.Range("J132").Value = 2640.6924
$ws.Range("K132").Value = 10347118.5
$ws.Range("L132").Value = 7922.0772
$ws.Range("M132").Value = -10344588.5
$ws.Range("N132").Value = -12982.0772
# Row 136
$ws.Range("H136").Value = 2656.8132
$ws.Range("I136").Value = 2820.164
$ws.Range("J136").Value = 1288.75
$ws.Range("K136").Value = 8460.492
$ws.Range("L136").Value = 3866.25
$ws.Range("M136").Value = -5910.492
$ws.Range("N136").Value = -8966.25

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 956.4
$ws.Range("I113").Value = 742.25
$ws.Range("K113").Value = 2226.75
$ws.Range("M113").Value = -56.75
# Row 132
$ws.Range("H132").Value = 298.62
$ws.Range("I132").Value = 222.81133
$ws.Range("J132").Value = 384.10638
$ws.Range("K132").Value = 668.43399
$ws.Range("L132").Value = 1152.31914
$ws.Range("M132").Value = 1861.56601
$ws.Range("N132").Value = -6212.31914
# Row 136
$ws.Range("H136").Value = 10326970
$ws.Range("I136").Value = 5251800.5
$ws.Range("K136").Value = 15755401.5
$ws.Range("M136").Value = -15752851.5

Write-Output "Updated 162 cells across 8 sheets"
